$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# -------------------------------------------------------------------------
# 1. Update the "VALOR MORA" (total overdue amount) and "Cant. Periodos"
#    (period count) summary figures at the top of the statement.
# -------------------------------------------------------------------------
$ws.Range("E11").Value = 428948
$ws.Range("F13").Value = 8

# -------------------------------------------------------------------------
# 2. Make room for a new payment-history row. The existing table's last
#    data row is row 22; insert a blank row right after it (before the old
#    row 23) so the closing "firma" block (old rows 27-28) shifts down to
#    rows 28-29, matching a normal Excel row insert.
# -------------------------------------------------------------------------
$ws.Rows("23").Insert()

# Carry the bottom-border "closing" formatting that used to belong to the
# last table row (old row 22) down onto the freshly inserted row 23 ...
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)  # xlPasteFormats

# ... and give the now-interior row 22 the regular (non-closing) row
# formatting used by every other data row, copied from row 21.
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# -------------------------------------------------------------------------
# 3. Re-populate the "Periodo Mora" / "Valor Mora" columns for the table
#    (rows 16-23) in ascending period order 2501..2508. Period 2501 keeps
#    its distinctive overdue value (30368); the new period 2508 row uses
#    the standard monthly value (56940).
# -------------------------------------------------------------------------
$periods = @(16,17,18,19,20,21,22,23)
$values  = @(30368,56940,56940,56940,56940,56940,56940,56940)
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $periods[$i]
    $period = 2501 + $i
    $ws.Cells.Item($row, 5).Value = [string]$period   # column E - Periodo Mora
    $ws.Cells.Item($row, 6).Value = $values[$i]        # column F - Valor Mora
}

# New row 23 needs the rest of its worker/detail data filled in too (same
# worker repeated for every period row, as in the rest of the table).
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1193319425"
$ws.Range("D23").Value = "EDILBERTO CASTRO TORRES"
$ws.Range("G23").Value = 1423500
